$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.510453462600708
$ws.Range("B1").Value = 1.65189802646637
$ws.Range("C1").Value = 3.788183450698853
$ws.Range("D1").Value = 2.294786691665649
$ws.Range("E1").Value = 0.8325689435005188
